$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5770
$ws1.Range("F6").Value = 617
$ws1.Range("F9").Value = 4496
$ws1.Range("F10").Value = 1807
$ws1.Range("F13").Value = 3035
$ws1.Range("F15").Value = 577
$ws1.Range("F16").Value = 234
$ws1.Range("F17").Value = 575
$ws1.Range("F18").Value = 495
$ws1.Range("F22").Value = 1743
$ws1.Range("F23").Value = 1277
$ws1.Range("F25").Value = 1501
$ws1.Range("F26").Value = 121
$ws1.Range("F29").Value = 522
$ws1.Range("F31").Value = 33
$ws1.Range("F32").Value = 77
$ws1.Range("F33").Value = 115
$ws1.Range("F35").Value = 3297
$ws1.Range("F36").Value = 731
$ws1.Range("F37").Value = 56
$ws1.Range("F38").Value = 210
$ws1.Range("F40").Value = 1559

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5770
$ws4.Range("F6").Value = 617
$ws4.Range("F9").Value = 4496
$ws4.Range("F10").Value = 1807
$ws4.Range("F13").Value = 3035
$ws4.Range("F15").Value = 577
$ws4.Range("F16").Value = 234
$ws4.Range("F17").Value = 575
$ws4.Range("F18").Value = 495
$ws4.Range("F23").Value = 1743
$ws4.Range("F24").Value = 1277
$ws4.Range("F26").Value = 1501
$ws4.Range("F27").Value = 121
$ws4.Range("F30").Value = 522
$ws4.Range("F32").Value = 33
$ws4.Range("F33").Value = 77
$ws4.Range("F34").Value = 115
$ws4.Range("F36").Value = 3297
$ws4.Range("F38").Value = 731
$ws4.Range("F39").Value = 56
$ws4.Range("F40").Value = 210
$ws4.Range("F42").Value = 1559
